$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 728-753 with new data (per diff) ---

$ws.Cells.Item(728, 4).Value = 45147
$ws.Cells.Item(728, 11).Value = "Lane Late"
$ws.Cells.Item(728, 12).Value = "Primera"
$ws.Cells.Item(728, 13).Value = 60
$ws.Cells.Item(728, 14).Value = 9000
$ws.Cells.Item(728, 15).Value = 9000
$ws.Cells.Item(728, 16).Value = 9000
$ws.Cells.Item(728, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(728, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(728, 19).Value = 600

$ws.Cells.Item(729, 4).Value = 45147
$ws.Cells.Item(729, 11).Value = "Navel Late"
$ws.Cells.Item(729, 12).Value = "Especial"
$ws.Cells.Item(729, 13).Value = 60
$ws.Cells.Item(729, 14).Value = 10000
$ws.Cells.Item(729, 15).Value = 10000
$ws.Cells.Item(729, 16).Value = 10000
$ws.Cells.Item(729, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(729, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(729, 19).Value = 667

$ws.Cells.Item(730, 4).Value = 45147
$ws.Cells.Item(730, 11).Value = "Navel Late"
$ws.Cells.Item(730, 12).Value = "Primera"
$ws.Cells.Item(730, 13).Value = 120
$ws.Cells.Item(730, 14).Value = 8000
$ws.Cells.Item(730, 15).Value = 8000
$ws.Cells.Item(730, 16).Value = 8000
$ws.Cells.Item(730, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(730, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(730, 19).Value = 533

$ws.Cells.Item(731, 4).Value = 45147
$ws.Cells.Item(731, 11).Value = "Navel Late"
$ws.Cells.Item(731, 12).Value = "Segunda"
$ws.Cells.Item(731, 13).Value = 100
$ws.Cells.Item(731, 14).Value = 6000
$ws.Cells.Item(731, 15).Value = 6000
$ws.Cells.Item(731, 16).Value = 6000
$ws.Cells.Item(731, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(731, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(731, 19).Value = 400

$ws.Cells.Item(732, 4).Value = 44756
$ws.Cells.Item(732, 11).Value = "Fukumoto"
$ws.Cells.Item(732, 12).Value = "Primera"
$ws.Cells.Item(732, 13).Value = 160
$ws.Cells.Item(732, 14).Value = 7500
$ws.Cells.Item(732, 15).Value = 8000
$ws.Cells.Item(732, 16).Value = 7750
$ws.Cells.Item(732, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(732, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(732, 19).Value = 517

$ws.Cells.Item(733, 4).Value = 44756
$ws.Cells.Item(733, 11).Value = "Fukumoto"
$ws.Cells.Item(733, 12).Value = "Segunda"
$ws.Cells.Item(733, 13).Value = 80
$ws.Cells.Item(733, 14).Value = 7000
$ws.Cells.Item(733, 15).Value = 7000
$ws.Cells.Item(733, 16).Value = 7000
$ws.Cells.Item(733, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(733, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(733, 19).Value = 467

$ws.Cells.Item(734, 4).Value = 44579
$ws.Cells.Item(734, 11).Value = "Valencia"
$ws.Cells.Item(734, 12).Value = "Primera"
$ws.Cells.Item(734, 13).Value = 120
$ws.Cells.Item(734, 14).Value = 9000
$ws.Cells.Item(734, 15).Value = 10000
$ws.Cells.Item(734, 16).Value = 9500
$ws.Cells.Item(734, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(734, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(734, 19).Value = 633

$ws.Cells.Item(735, 4).Value = 44343
$ws.Cells.Item(735, 11).Value = "Fukumoto"
$ws.Cells.Item(735, 12).Value = "Primera"
$ws.Cells.Item(735, 13).Value = 120
$ws.Cells.Item(735, 14).Value = 12000
$ws.Cells.Item(735, 15).Value = 13000
$ws.Cells.Item(735, 16).Value = 12500
$ws.Cells.Item(735, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(735, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(735, 19).Value = 833

$ws.Cells.Item(736, 4).Value = 44343
$ws.Cells.Item(736, 11).Value = "Fukumoto"
$ws.Cells.Item(736, 12).Value = "Segunda"
$ws.Cells.Item(736, 13).Value = 60
$ws.Cells.Item(736, 14).Value = 11000
$ws.Cells.Item(736, 15).Value = 11000
$ws.Cells.Item(736, 16).Value = 11000
$ws.Cells.Item(736, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(736, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(736, 19).Value = 733

$ws.Cells.Item(737, 4).Value = 44526
$ws.Cells.Item(737, 11).Value = "Lane Late"
$ws.Cells.Item(737, 12).Value = "Primera"
$ws.Cells.Item(737, 13).Value = 160
$ws.Cells.Item(737, 14).Value = 8500
$ws.Cells.Item(737, 15).Value = 9000
$ws.Cells.Item(737, 16).Value = 8750
$ws.Cells.Item(737, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(737, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(737, 19).Value = 583

$ws.Cells.Item(738, 4).Value = 44341
$ws.Cells.Item(738, 11).Value = "Fukumoto"
$ws.Cells.Item(738, 12).Value = "Primera"
$ws.Cells.Item(738, 13).Value = 120
$ws.Cells.Item(738, 14).Value = 12000
$ws.Cells.Item(738, 15).Value = 13000
$ws.Cells.Item(738, 16).Value = 12500
$ws.Cells.Item(738, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(738, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(738, 19).Value = 833

$ws.Cells.Item(739, 4).Value = 44341
$ws.Cells.Item(739, 11).Value = "Fukumoto"
$ws.Cells.Item(739, 12).Value = "Segunda"
$ws.Cells.Item(739, 13).Value = 60
$ws.Cells.Item(739, 14).Value = 11000
$ws.Cells.Item(739, 15).Value = 11000
$ws.Cells.Item(739, 16).Value = 11000
$ws.Cells.Item(739, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(739, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(739, 19).Value = 733

$ws.Cells.Item(740, 4).Value = 44397
$ws.Cells.Item(740, 11).Value = "Fukumoto"
$ws.Cells.Item(740, 12).Value = "Primera"
$ws.Cells.Item(740, 13).Value = 120
$ws.Cells.Item(740, 14).Value = 7500
$ws.Cells.Item(740, 15).Value = 8000
$ws.Cells.Item(740, 16).Value = 7750
$ws.Cells.Item(740, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(740, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(740, 19).Value = 517

$ws.Cells.Item(741, 4).Value = 44397
$ws.Cells.Item(741, 11).Value = "Fukumoto"
$ws.Cells.Item(741, 12).Value = "Segunda"
$ws.Cells.Item(741, 13).Value = 80
$ws.Cells.Item(741, 14).Value = 6500
$ws.Cells.Item(741, 15).Value = 6500
$ws.Cells.Item(741, 16).Value = 6500
$ws.Cells.Item(741, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(741, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(741, 19).Value = 433

$ws.Cells.Item(742, 4).Value = 45068
$ws.Cells.Item(742, 11).Value = "Valencia"
$ws.Cells.Item(742, 12).Value = "Especial"
$ws.Cells.Item(742, 13).Value = 100
$ws.Cells.Item(742, 14).Value = 14000
$ws.Cells.Item(742, 15).Value = 14000
$ws.Cells.Item(742, 16).Value = 14000
$ws.Cells.Item(742, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(742, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(742, 19).Value = 933

$ws.Cells.Item(743, 4).Value = 45068
$ws.Cells.Item(743, 11).Value = "Valencia"
$ws.Cells.Item(743, 12).Value = "Primera"
$ws.Cells.Item(743, 13).Value = 120
$ws.Cells.Item(743, 14).Value = 12000
$ws.Cells.Item(743, 15).Value = 12000
$ws.Cells.Item(743, 16).Value = 12000
$ws.Cells.Item(743, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(743, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(743, 19).Value = 800

$ws.Cells.Item(744, 4).Value = 44631
$ws.Cells.Item(744, 11).Value = "Valencia"
$ws.Cells.Item(744, 12).Value = "Primera"
$ws.Cells.Item(744, 13).Value = 120
$ws.Cells.Item(744, 14).Value = 10000
$ws.Cells.Item(744, 15).Value = 11000
$ws.Cells.Item(744, 16).Value = 10500
$ws.Cells.Item(744, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(744, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(744, 19).Value = 700

$ws.Cells.Item(745, 4).Value = 44631
$ws.Cells.Item(745, 11).Value = "Valencia"
$ws.Cells.Item(745, 12).Value = "Segunda"
$ws.Cells.Item(745, 13).Value = 60
$ws.Cells.Item(745, 14).Value = 9000
$ws.Cells.Item(745, 15).Value = 9000
$ws.Cells.Item(745, 16).Value = 9000
$ws.Cells.Item(745, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(745, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(745, 19).Value = 600

$ws.Cells.Item(746, 4).Value = 44645
$ws.Cells.Item(746, 11).Value = "Valencia"
$ws.Cells.Item(746, 12).Value = "Primera"
$ws.Cells.Item(746, 13).Value = 120
$ws.Cells.Item(746, 14).Value = 10000
$ws.Cells.Item(746, 15).Value = 11000
$ws.Cells.Item(746, 16).Value = 10500
$ws.Cells.Item(746, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(746, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(746, 19).Value = 700

$ws.Cells.Item(747, 4).Value = 44645
$ws.Cells.Item(747, 11).Value = "Valencia"
$ws.Cells.Item(747, 12).Value = "Segunda"
$ws.Cells.Item(747, 13).Value = 120
$ws.Cells.Item(747, 14).Value = 8000
$ws.Cells.Item(747, 15).Value = 9000
$ws.Cells.Item(747, 16).Value = 8500
$ws.Cells.Item(747, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(747, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(747, 19).Value = 567

$ws.Cells.Item(748, 4).Value = 44953
$ws.Cells.Item(748, 11).Value = "Valencia"
$ws.Cells.Item(748, 12).Value = "Primera"
$ws.Cells.Item(748, 13).Value = 60
$ws.Cells.Item(748, 14).Value = 10000
$ws.Cells.Item(748, 15).Value = 10000
$ws.Cells.Item(748, 16).Value = 10000
$ws.Cells.Item(748, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(748, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(748, 19).Value = 667

$ws.Cells.Item(749, 4).Value = 44953
$ws.Cells.Item(749, 11).Value = "Valencia"
$ws.Cells.Item(749, 12).Value = "Segunda"
$ws.Cells.Item(749, 13).Value = 60
$ws.Cells.Item(749, 14).Value = 8000
$ws.Cells.Item(749, 15).Value = 8000
$ws.Cells.Item(749, 16).Value = 8000
$ws.Cells.Item(749, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(749, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(749, 19).Value = 533

$ws.Cells.Item(750, 4).Value = 44417
$ws.Cells.Item(750, 11).Value = "Fukumoto"
$ws.Cells.Item(750, 12).Value = "Primera"
$ws.Cells.Item(750, 13).Value = 160
$ws.Cells.Item(750, 14).Value = 6000
$ws.Cells.Item(750, 15).Value = 6200
$ws.Cells.Item(750, 16).Value = 6100
$ws.Cells.Item(750, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(750, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(750, 19).Value = 407

$ws.Cells.Item(751, 4).Value = 44417
$ws.Cells.Item(751, 11).Value = "Fukumoto"
$ws.Cells.Item(751, 12).Value = "Segunda"
$ws.Cells.Item(751, 13).Value = 80
$ws.Cells.Item(751, 14).Value = 5500
$ws.Cells.Item(751, 15).Value = 5800
$ws.Cells.Item(751, 16).Value = 5650
$ws.Cells.Item(751, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(751, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(751, 19).Value = 377

$ws.Cells.Item(752, 4).Value = 44417
$ws.Cells.Item(752, 11).Value = "Navel Late"
$ws.Cells.Item(752, 12).Value = "Primera"
$ws.Cells.Item(752, 13).Value = 240
$ws.Cells.Item(752, 14).Value = 6000
$ws.Cells.Item(752, 15).Value = 6500
$ws.Cells.Item(752, 16).Value = 6250
$ws.Cells.Item(752, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(752, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(752, 19).Value = 417

$ws.Cells.Item(753, 4).Value = 44417
$ws.Cells.Item(753, 11).Value = "Navel Late"
$ws.Cells.Item(753, 12).Value = "Segunda"
$ws.Cells.Item(753, 13).Value = 120
$ws.Cells.Item(753, 14).Value = 5000
$ws.Cells.Item(753, 15).Value = 5500
$ws.Cells.Item(753, 16).Value = 5250
$ws.Cells.Item(753, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(753, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(753, 19).Value = 350

# --- Append new rows 754-758 (full rows, constant columns + varying data) ---

$ws.Cells.Item(754, 1).Value = 7
$ws.Cells.Item(754, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(754, 3).Value = "Ñuble"
$ws.Cells.Item(754, 4).Value = 44313
$ws.Cells.Item(754, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(754, 5).Value = 16
$ws.Cells.Item(754, 6).Value = "Fruta"
$ws.Cells.Item(754, 7).Value = 100102
$ws.Cells.Item(754, 8).Value = "Cítricos"
$ws.Cells.Item(754, 9).Value = 100102005
$ws.Cells.Item(754, 10).Value = "Naranja"
$ws.Cells.Item(754, 11).Value = "Fukumoto"
$ws.Cells.Item(754, 12).Value = "Primera"
$ws.Cells.Item(754, 13).Value = 160
$ws.Cells.Item(754, 14).Value = 14000
$ws.Cells.Item(754, 15).Value = 15000
$ws.Cells.Item(754, 16).Value = 14500
$ws.Cells.Item(754, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(754, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(754, 19).Value = 967
$ws.Cells.Item(754, 20).Value = 15

$ws.Cells.Item(755, 1).Value = 7
$ws.Cells.Item(755, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(755, 3).Value = "Ñuble"
$ws.Cells.Item(755, 4).Value = 44313
$ws.Cells.Item(755, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(755, 5).Value = 16
$ws.Cells.Item(755, 6).Value = "Fruta"
$ws.Cells.Item(755, 7).Value = 100102
$ws.Cells.Item(755, 8).Value = "Cítricos"
$ws.Cells.Item(755, 9).Value = 100102005
$ws.Cells.Item(755, 10).Value = "Naranja"
$ws.Cells.Item(755, 11).Value = "Fukumoto"
$ws.Cells.Item(755, 12).Value = "Segunda"
$ws.Cells.Item(755, 13).Value = 80
$ws.Cells.Item(755, 14).Value = 13000
$ws.Cells.Item(755, 15).Value = 13000
$ws.Cells.Item(755, 16).Value = 13000
$ws.Cells.Item(755, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(755, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(755, 19).Value = 867
$ws.Cells.Item(755, 20).Value = 15

$ws.Cells.Item(756, 1).Value = 7
$ws.Cells.Item(756, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(756, 3).Value = "Ñuble"
$ws.Cells.Item(756, 4).Value = 44979
$ws.Cells.Item(756, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(756, 5).Value = 16
$ws.Cells.Item(756, 6).Value = "Fruta"
$ws.Cells.Item(756, 7).Value = 100102
$ws.Cells.Item(756, 8).Value = "Cítricos"
$ws.Cells.Item(756, 9).Value = 100102005
$ws.Cells.Item(756, 10).Value = "Naranja"
$ws.Cells.Item(756, 11).Value = "Valencia"
$ws.Cells.Item(756, 12).Value = "Primera"
$ws.Cells.Item(756, 13).Value = 400
$ws.Cells.Item(756, 14).Value = 13000
$ws.Cells.Item(756, 15).Value = 14000
$ws.Cells.Item(756, 16).Value = 13500
$ws.Cells.Item(756, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(756, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(756, 19).Value = 900
$ws.Cells.Item(756, 20).Value = 15

$ws.Cells.Item(757, 1).Value = 7
$ws.Cells.Item(757, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(757, 3).Value = "Ñuble"
$ws.Cells.Item(757, 4).Value = 44364
$ws.Cells.Item(757, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(757, 5).Value = 16
$ws.Cells.Item(757, 6).Value = "Fruta"
$ws.Cells.Item(757, 7).Value = 100102
$ws.Cells.Item(757, 8).Value = "Cítricos"
$ws.Cells.Item(757, 9).Value = 100102005
$ws.Cells.Item(757, 10).Value = "Naranja"
$ws.Cells.Item(757, 11).Value = "Fukumoto"
$ws.Cells.Item(757, 12).Value = "Primera"
$ws.Cells.Item(757, 13).Value = 120
$ws.Cells.Item(757, 14).Value = 9000
$ws.Cells.Item(757, 15).Value = 9500
$ws.Cells.Item(757, 16).Value = 9250
$ws.Cells.Item(757, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(757, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(757, 19).Value = 617
$ws.Cells.Item(757, 20).Value = 15

$ws.Cells.Item(758, 1).Value = 7
$ws.Cells.Item(758, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(758, 3).Value = "Ñuble"
$ws.Cells.Item(758, 4).Value = 44364
$ws.Cells.Item(758, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(758, 5).Value = 16
$ws.Cells.Item(758, 6).Value = "Fruta"
$ws.Cells.Item(758, 7).Value = 100102
$ws.Cells.Item(758, 8).Value = "Cítricos"
$ws.Cells.Item(758, 9).Value = 100102005
$ws.Cells.Item(758, 10).Value = "Naranja"
$ws.Cells.Item(758, 11).Value = "Fukumoto"
$ws.Cells.Item(758, 12).Value = "Segunda"
$ws.Cells.Item(758, 13).Value = 120
$ws.Cells.Item(758, 14).Value = 8000
$ws.Cells.Item(758, 15).Value = 8500
$ws.Cells.Item(758, 16).Value = 8250
$ws.Cells.Item(758, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(758, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(758, 19).Value = 550
$ws.Cells.Item(758, 20).Value = 15

